# Apply updated cryptos data (prices and 1h volume deltas) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.272.88'
$ws.Range("E2").Value = '  -4.69%  '
$ws.Range("D3").Value = '3.259.11'
$ws.Range("E3").Value = '  -7.26%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''598.37'
$ws.Range("E5").Value = '  -3.23%  '
$ws.Range("D6").Value = '''151.15'
$ws.Range("E6").Value = '  -12.47%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.252.15'
$ws.Range("E8").Value = '  -7.40%  '
$ws.Range("E9").Value = '  -11.18%  '
$ws.Range("E10").Value = '  -12.84%  '
$ws.Range("D11").Value = '''6.74'
$ws.Range("E11").Value = '  -4.71%  '
$ws.Range("D12").Value = '''0.506'
$ws.Range("E12").Value = '  -13.63%  '
$ws.Range("D13").Value = '''38.07'
$ws.Range("E13").Value = '  -17.80%  '
$ws.Range("E14").Value = '  -11.64%  '
$ws.Range("D15").Value = '3.783.41'
$ws.Range("E15").Value = '  -7.53%  '
$ws.Range("D16").Value = '67.321.23'
$ws.Range("E16").Value = '  -4.87%  '
$ws.Range("D17").Value = '3.260.09'
$ws.Range("E17").Value = '  -7.47%  '
$ws.Range("D18").Value = '''544.35'
$ws.Range("E18").Value = '  -10.33%  '
$ws.Range("E19").Value = '  -6.02%  '
$ws.Range("E20").Value = '  -13.62%  '
$ws.Range("E21").Value = '  -14.43%  '
$ws.Range("D22").Value = '''0.762'
$ws.Range("E22").Value = '  -13.36%  '
$ws.Range("D23").Value = '''7.85'
$ws.Range("E23").Value = '  -14.41%  '
$ws.Range("D24").Value = '''85.47'
$ws.Range("E24").Value = '  -12.69%  '
$ws.Range("E25").Value = '  -12.94%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '''3.26'
$ws.Range("E27").Value = '  -12.44%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''8.06'
$ws.Range("E28").Value = '  -10.64%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '''29.36'
$ws.Range("E29").Value = '  -12.55%  '
$ws.Range("D30").Value = '''2.13'
$ws.Range("E30").Value = '  -17.05%  '
$ws.Range("E31").Value = '  -10.61%  '
$ws.Range("E32").Value = '  -11.79%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '''546.74'
$ws.Range("E33").Value = '  -14.41%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''6.65'
$ws.Range("E34").Value = '  -17.48%  '
$ws.Range("D35").Value = '''5.71'
$ws.Range("E35").Value = '  -15.76%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '''0.0447'
$ws.Range("E37").Value = '  -7.65%  '
$ws.Range("D38").Value = '''53.53'
$ws.Range("E38").Value = '  -5.38%  '
$ws.Range("D39").Value = '''0.0854'
$ws.Range("E39").Value = '  -14.18%  '
$ws.Range("D40").Value = '''9.17'
$ws.Range("E40").Value = '  -14.94%  '
$ws.Range("D41").Value = '''0.128'
$ws.Range("E41").Value = '  -9.68%  '
$ws.Range("D42").Value = '2.929.56'
$ws.Range("E42").Value = '  -12.41%  '
$ws.Range("E43").Value = '  -21.80%  '
$ws.Range("E44").Value = '  -15.71%  '
$ws.Range("D45").Value = '0.0₃0584'
$ws.Range("E45").Value = '  -18.07%  '
$ws.Range("D46").Value = '''2.19'
$ws.Range("E46").Value = '  -13.17%  '
$ws.Range("D47").Value = '''26.46'
$ws.Range("E47").Value = '  -16.48%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '''128.07'
$ws.Range("E49").Value = '  -4.82%  '
$ws.Range("D50").Value = '''2.36'
$ws.Range("E50").Value = '  -19.47%  '
$ws.Range("E51").Value = '  -12.59%  '
